$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9421.857
$ws.Range("I40").Value = 5190.6
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 5190.6
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = -5015.6
$ws.Range("N40").Value = -20350
$ws.Range("H62").Value = 3891.5557
$ws.Range("I62").Value = 3589.1428
$ws.Range("J62").Value = 4950
$ws.Range("K62").Value = 3589.1428
$ws.Range("L62").Value = 4950
$ws.Range("M62").Value = -2965.1428
$ws.Range("N62").Value = -6198
$ws.Range("H65").Value = 3891.5557
$ws.Range("I65").Value = 3589.1428
$ws.Range("J65").Value = 4950
$ws.Range("K65").Value = 17945.714
$ws.Range("L65").Value = 24750
$ws.Range("M65").Value = -14825.714
$ws.Range("N65").Value = -30990
$ws.Range("H70").Value = 80645.30499999999
$ws.Range("I70").Value = 3231.6667
$ws.Range("J70").Value = 146999.86
$ws.Range("K70").Value = 9695.000100000001
$ws.Range("L70").Value = 440999.58
$ws.Range("M70").Value = -9425.000100000001
$ws.Range("N70").Value = -441539.58
$ws.Range("H73").Value = 80645.30499999999
$ws.Range("I73").Value = 3231.6667
$ws.Range("J73").Value = 146999.86
$ws.Range("K73").Value = 9695.000100000001
$ws.Range("L73").Value = 440999.58
$ws.Range("M73").Value = -8759.000100000001
$ws.Range("N73").Value = -442871.58
$ws.Range("H80").Value = 1195.5
$ws.Range("I80").Value = 697
$ws.Range("J80").Value = 1361.6666
$ws.Range("K80").Value = 2091
$ws.Range("L80").Value = 4084.9998
$ws.Range("M80").Value = -1093
$ws.Range("N80").Value = -6080.9998
$ws.Range("H82").Value = 2695.9092
$ws.Range("I82").Value = 1565.5
$ws.Range("K82").Value = 4696.5
$ws.Range("M82").Value = -4290.5
$ws.Range("H83").Value = 1195.5
$ws.Range("I83").Value = 697
$ws.Range("J83").Value = 1361.6666
$ws.Range("K83").Value = 6273
$ws.Range("L83").Value = 12254.9994
$ws.Range("M83").Value = -1281
$ws.Range("N83").Value = -22238.9994
$ws.Range("H85").Value = 2695.9092
$ws.Range("I85").Value = 1565.5
$ws.Range("K85").Value = 4696.5
$ws.Range("M85").Value = -3292.5
$ws.Range("H86").Value = 5432.2
$ws.Range("J86").Value = 5139.375
$ws.Range("L86").Value = 5139.375
$ws.Range("N86").Value = -7385.375
$ws.Range("H89").Value = 5432.2
$ws.Range("J89").Value = 5139.375
$ws.Range("L89").Value = 25696.875
$ws.Range("N89").Value = -36928.875
$ws.Range("H107").Value = 1889.5883
$ws.Range("I107").Value = 1741.6666
$ws.Range("J107").Value = 2999
$ws.Range("K107").Value = 1741.6666
$ws.Range("L107").Value = 2999
$ws.Range("M107").Value = 178.3334
$ws.Range("N107").Value = -6839
$ws.Range("H112").Value = 4527.385
$ws.Range("J112").Value = 4779.6665
$ws.Range("L112").Value = 14338.9995
$ws.Range("N112").Value = -16554.9995
$ws.Range("H113").Value = 6089.5
$ws.Range("I113").Value = 9232.333000000001
$ws.Range("J113").Value = 5041.8887
$ws.Range("K113").Value = 9232.333000000001
$ws.Range("L113").Value = 5041.8887
$ws.Range("M113").Value = -5978.333000000001
$ws.Range("N113").Value = -11549.8887
$ws.Range("H115").Value = 1244.7778
$ws.Range("I115").Value = 1244.7778
$ws.Range("K115").Value = 3734.3334
$ws.Range("M115").Value = -2167.3334
$ws.Range("H129").Value = 914.55554
$ws.Range("I129").Value = 914.55554
$ws.Range("K129").Value = 2743.66662
$ws.Range("M129").Value = 2256.33338
$ws.Range("H134").Value = 94999
$ws.Range("J134").Value = 94999
$ws.Range("L134").Value = 94999
$ws.Range("N134").Value = -105139
$ws.Range("H136").Value = 67998.2
$ws.Range("J136").Value = 67998.2
$ws.Range("L136").Value = 67998.2
$ws.Range("N136").Value = -78198.2
$ws.Range("H140").Value = 94999
$ws.Range("J140").Value = 94999
$ws.Range("L140").Value = 94999
$ws.Range("N140").Value = -105359

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15774.941
$ws.Range("I45").Value = 24741.555
$ws.Range("K45").Value = 24741.555
$ws.Range("M45").Value = -24364.555
$ws.Range("H46").Value = 29576
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H97").Value = 2040.72
$ws.Range("I97").Value = 927.26666
$ws.Range("J97").Value = 3710.9
$ws.Range("K97").Value = 927.26666
$ws.Range("L97").Value = 3710.9
$ws.Range("M97").Value = -431.26666
$ws.Range("N97").Value = -4702.9
$ws.Range("H132").Value = 2540.1667
$ws.Range("I132").Value = 2625.6365
$ws.Range("K132").Value = 7876.9095
$ws.Range("M132").Value = -5346.9095
$ws.Range("H141").Value = 89429
$ws.Range("J141").Value = 89429
$ws.Range("L141").Value = 89429
$ws.Range("N141").Value = -99789

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 43748.75
$ws.Range("J133").Value = 43748.75
$ws.Range("L133").Value = 43748.75
$ws.Range("N133").Value = -53868.75
$ws.Range("H134").Value = 3244.4443
$ws.Range("I134").Value = 3244.4443
$ws.Range("K134").Value = 9733.332900000001
$ws.Range("M134").Value = -7198.332900000001
$ws.Range("H139").Value = 77139
$ws.Range("J139").Value = 80708
$ws.Range("L139").Value = 80708
$ws.Range("N139").Value = -90988

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10639.0625
$ws.Range("I31").Value = 18218.625
$ws.Range("K31").Value = 18218.625
$ws.Range("M31").Value = -17923.625
$ws.Range("H34").Value = 10639.0625
$ws.Range("I34").Value = 18218.625
$ws.Range("K34").Value = 18218.625
$ws.Range("M34").Value = -18016.625
$ws.Range("H99").Value = 4999
$ws.Range("I99").Value = 5694.75
$ws.Range("K99").Value = 5694.75
$ws.Range("M99").Value = -4196.75
$ws.Range("H105").Value = 1161.5454
$ws.Range("I105").Value = 1197.25
$ws.Range("K105").Value = 1197.25
$ws.Range("M105").Value = 549.75
$ws.Range("H126").Value = 4999
$ws.Range("I126").Value = 5694.75
$ws.Range("K126").Value = 17084.25
$ws.Range("M126").Value = -14614.25
$ws.Range("H137").Value = 69419.664
$ws.Range("J137").Value = 69999.625
$ws.Range("L137").Value = 69999.625
$ws.Range("N137").Value = -80199.625

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 143.8
$ws.Range("I40").Value = 73
$ws.Range("K40").Value = 292
$ws.Range("M40").Value = -223
$ws.Range("H51").Value = 2999.7144
$ws.Range("J51").Value = 3333
$ws.Range("L51").Value = 9999
$ws.Range("N51").Value = -10919
$ws.Range("H107").Value = 737.6087
$ws.Range("I107").Value = 290.83334
$ws.Range("K107").Value = 872.5000200000001
$ws.Range("M107").Value = 1047.49998
$ws.Range("H131").Value = 2273.3057
$ws.Range("J131").Value = 2396.6562
$ws.Range("L131").Value = 7189.9686
$ws.Range("N131").Value = -17269.9686

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 187316
$ws.Range("J121").Value = 187316
$ws.Range("L121").Value = 187316
$ws.Range("N121").Value = -190810
$ws.Range("H132").Value = 7998.5
$ws.Range("J132").Value = 8998.5
$ws.Range("L132").Value = 26995.5
$ws.Range("N132").Value = -32055.5
$ws.Range("H136").Value = 31808.555
$ws.Range("J136").Value = 31808.555
$ws.Range("L136").Value = 95425.66500000001
$ws.Range("N136").Value = -100525.665
$ws.Range("H137").Value = 88999
$ws.Range("J137").Value = 88999
$ws.Range("L137").Value = 88999
$ws.Range("N137").Value = -99199

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1791
$ws.Range("I22").Value = 1576.3
$ws.Range("J22").Value = 2220.4
$ws.Range("K22").Value = 1576.3
$ws.Range("L22").Value = 2220.4
$ws.Range("M22").Value = -1281.3
$ws.Range("N22").Value = -2810.4
$ws.Range("H27").Value = 1791
$ws.Range("I27").Value = 1576.3
$ws.Range("J27").Value = 2220.4
$ws.Range("K27").Value = 1576.3
$ws.Range("L27").Value = 2220.4
$ws.Range("M27").Value = -1469.3
$ws.Range("N27").Value = -2434.4
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1505.5714
$ws.Range("I100").Value = 1505.5714
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3011.1428
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2470.1428
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 2402.0417
$ws.Range("I107").Value = 2383.2
$ws.Range("K107").Value = 7149.599999999999
$ws.Range("M107").Value = -5229.599999999999
$ws.Range("H132").Value = 3664.3635
$ws.Range("I132").Value = 3559.8
$ws.Range("K132").Value = 10679.4
$ws.Range("M132").Value = -8149.400000000001
$ws.Range("H136").Value = 1440.6364
$ws.Range("I136").Value = 1284.8823
$ws.Range("K136").Value = 3854.6469
$ws.Range("M136").Value = -1304.6469
$ws.Range("H138").Value = 96493.25
$ws.Range("J138").Value = 98657.664
$ws.Range("L138").Value = 98657.664
$ws.Range("N138").Value = -108937.664
